# Auto-generated script applying cell-value updates per the commit diff.
# Values are the authoritative computed/cached numbers for the affected Leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 822.76746   # ALC!H17
$ws.Cells.Item(17, 10).Value = 822.76746   # ALC!J17
$ws.Cells.Item(17, 12).Value = 2468.30238   # ALC!L17
$ws.Cells.Item(17, 14).Value = -2804.30238   # ALC!N17
$ws.Cells.Item(92, 8).Value = 368.20834   # ALC!H92
$ws.Cells.Item(92, 9).Value = 297.35   # ALC!I92
$ws.Cells.Item(92, 10).Value = 722.5   # ALC!J92
$ws.Cells.Item(92, 11).Value = 297.35   # ALC!K92
$ws.Cells.Item(92, 12).Value = 722.5   # ALC!L92
$ws.Cells.Item(92, 13).Value = 950.65   # ALC!M92
$ws.Cells.Item(92, 14).Value = -3218.5   # ALC!N92
$ws.Cells.Item(101, 8).Value = 593.9231   # ALC!H101
$ws.Cells.Item(101, 9).Value = 494.16666   # ALC!I101
$ws.Cells.Item(101, 10).Value = 679.4286   # ALC!J101
$ws.Cells.Item(101, 11).Value = 1482.49998   # ALC!K101
$ws.Cells.Item(101, 12).Value = 2038.2858   # ALC!L101
$ws.Cells.Item(101, 13).Value = 139.5000199999999   # ALC!M101
$ws.Cells.Item(101, 14).Value = -5282.2858   # ALC!N101
$ws.Cells.Item(103, 8).Value = 866.3333   # ALC!H103
$ws.Cells.Item(103, 10).Value = 959.6   # ALC!J103
$ws.Cells.Item(103, 12).Value = 2878.8   # ALC!L103
$ws.Cells.Item(103, 14).Value = -4050.8   # ALC!N103

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1105   # ARM!H102
$ws.Cells.Item(102, 9).Value = 980.7692   # ARM!I102
$ws.Cells.Item(102, 10).Value = 1643.3334   # ARM!J102
$ws.Cells.Item(102, 11).Value = 980.7692   # ARM!K102
$ws.Cells.Item(102, 12).Value = 1643.3334   # ARM!L102
$ws.Cells.Item(102, 13).Value = 641.2308   # ARM!M102
$ws.Cells.Item(102, 14).Value = -4887.3334   # ARM!N102
$ws.Cells.Item(139, 8).Value = 201555.2   # ARM!H139
$ws.Cells.Item(139, 10).Value = 201555.2   # ARM!J139
$ws.Cells.Item(139, 12).Value = 201555.2   # ARM!L139
$ws.Cells.Item(139, 14).Value = -211835.2   # ARM!N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 500   # BSM!H15
$ws.Cells.Item(15, 9).Value = 500   # BSM!I15
$ws.Cells.Item(15, 11).Value = 500   # BSM!K15
$ws.Cells.Item(15, 13).Value = -273   # BSM!M15
$ws.Cells.Item(23, 8).Value = 114   # BSM!H23
$ws.Cells.Item(23, 9).Value = 0   # BSM!I23
$ws.Cells.Item(23, 10).Value = 114   # BSM!J23
$ws.Cells.Item(23, 11).Value = 0   # BSM!K23
$ws.Cells.Item(23, 12).Value = ""   # BSM!L23 (cleared)
$ws.Cells.Item(23, 13).Value = 114   # BSM!M23
$ws.Cells.Item(23, 14).Value = -680   # BSM!N23
$ws.Cells.Item(86, 8).Value = 3086.5   # BSM!H86
$ws.Cells.Item(86, 9).Value = 3283.6667   # BSM!I86
$ws.Cells.Item(86, 10).Value = 2968.2   # BSM!J86
$ws.Cells.Item(86, 11).Value = 3283.6667   # BSM!K86
$ws.Cells.Item(86, 12).Value = 2968.2   # BSM!L86
$ws.Cells.Item(86, 13).Value = -2160.6667   # BSM!M86
$ws.Cells.Item(86, 14).Value = -5214.2   # BSM!N86
$ws.Cells.Item(89, 8).Value = 3086.5   # BSM!H89
$ws.Cells.Item(89, 9).Value = 3283.6667   # BSM!I89
$ws.Cells.Item(89, 10).Value = 2968.2   # BSM!J89
$ws.Cells.Item(89, 11).Value = 16418.3335   # BSM!K89
$ws.Cells.Item(89, 12).Value = 14841   # BSM!L89
$ws.Cells.Item(89, 13).Value = -10802.3335   # BSM!M89
$ws.Cells.Item(89, 14).Value = -26073   # BSM!N89
$ws.Cells.Item(99, 8).Value = 83334700   # BSM!H99
$ws.Cells.Item(99, 9).Value = 142857860   # BSM!I99
$ws.Cells.Item(99, 11).Value = 142857860   # BSM!K99
$ws.Cells.Item(99, 13).Value = -142856362   # BSM!M99
$ws.Cells.Item(107, 8).Value = 1917.2   # BSM!H107
$ws.Cells.Item(107, 9).Value = 1787.5714   # BSM!I107
$ws.Cells.Item(107, 10).Value = 1987   # BSM!J107
$ws.Cells.Item(107, 11).Value = 1787.5714   # BSM!K107
$ws.Cells.Item(107, 12).Value = 1987   # BSM!L107
$ws.Cells.Item(107, 13).Value = 132.4286   # BSM!M107
$ws.Cells.Item(107, 14).Value = -5827   # BSM!N107
$ws.Cells.Item(134, 8).Value = 811.28766   # BSM!H134
$ws.Cells.Item(134, 9).Value = 745.07245   # BSM!I134
$ws.Cells.Item(134, 10).Value = 1953.5   # BSM!J134
$ws.Cells.Item(134, 11).Value = 2235.21735   # BSM!K134
$ws.Cells.Item(134, 12).Value = 5860.5   # BSM!L134
$ws.Cells.Item(134, 13).Value = 299.7826500000001   # BSM!M134
$ws.Cells.Item(134, 14).Value = -10930.5   # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 20030   # CRP!H74
$ws.Cells.Item(74, 9).Value = 90   # CRP!I74
$ws.Cells.Item(74, 10).Value = 30000   # CRP!J74
$ws.Cells.Item(74, 11).Value = 90   # CRP!K74
$ws.Cells.Item(74, 12).Value = 30000   # CRP!L74
$ws.Cells.Item(74, 13).Value = 784   # CRP!M74
$ws.Cells.Item(74, 14).Value = -31748   # CRP!N74
$ws.Cells.Item(77, 8).Value = 20030   # CRP!H77
$ws.Cells.Item(77, 9).Value = 90   # CRP!I77
$ws.Cells.Item(77, 10).Value = 30000   # CRP!J77
$ws.Cells.Item(77, 11).Value = 270   # CRP!K77
$ws.Cells.Item(77, 12).Value = 90000   # CRP!L77
$ws.Cells.Item(77, 13).Value = 4098   # CRP!M77
$ws.Cells.Item(77, 14).Value = -98736   # CRP!N77

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(95, 8).Value = 3000   # CUL!H95
$ws.Cells.Item(95, 10).Value = 3000   # CUL!J95
$ws.Cells.Item(95, 12).Value = 9000   # CUL!L95
$ws.Cells.Item(95, 14).Value = -13118   # CUL!N95
$ws.Cells.Item(96, 8).Value = 7000   # CUL!H96
$ws.Cells.Item(96, 10).Value = 7000   # CUL!J96
$ws.Cells.Item(96, 12).Value = 21000   # CUL!L96
$ws.Cells.Item(96, 14).Value = -25118   # CUL!N96
$ws.Cells.Item(99, 8).Value = 0   # CUL!H99
$ws.Cells.Item(99, 10).Value = 0   # CUL!J99
$ws.Cells.Item(99, 12).Value = ""   # CUL!L99 (cleared)
$ws.Cells.Item(99, 14).Value = 0   # CUL!N99
$ws.Cells.Item(100, 8).Value = 4007   # CUL!H100
$ws.Cells.Item(100, 10).Value = 4007   # CUL!J100
$ws.Cells.Item(100, 12).Value = 12021   # CUL!L100
$ws.Cells.Item(100, 14).Value = -13643   # CUL!N100
$ws.Cells.Item(102, 8).Value = 3000   # CUL!H102
$ws.Cells.Item(102, 9).Value = 3000   # CUL!I102
$ws.Cells.Item(102, 10).Value = 0   # CUL!J102
$ws.Cells.Item(102, 11).Value = 9000   # CUL!K102
$ws.Cells.Item(102, 12).Value = ""   # CUL!L102 (cleared)
$ws.Cells.Item(102, 13).Value = -6566   # CUL!M102
$ws.Cells.Item(102, 14).Value = 0   # CUL!N102
$ws.Cells.Item(104, 8).Value = 0   # CUL!H104
$ws.Cells.Item(104, 10).Value = 0   # CUL!J104
$ws.Cells.Item(104, 12).Value = ""   # CUL!L104 (cleared)
$ws.Cells.Item(104, 14).Value = 0   # CUL!N104
$ws.Cells.Item(108, 8).Value = 1363.5   # CUL!H108
$ws.Cells.Item(108, 9).Value = 1363.5   # CUL!I108
$ws.Cells.Item(108, 11).Value = 4090.5   # CUL!K108
$ws.Cells.Item(108, 13).Value = -1210.5   # CUL!M108
$ws.Cells.Item(109, 8).Value = 5000   # CUL!H109
$ws.Cells.Item(109, 9).Value = 0   # CUL!I109
$ws.Cells.Item(109, 10).Value = 5000   # CUL!J109
$ws.Cells.Item(109, 11).Value = 0   # CUL!K109
$ws.Cells.Item(109, 12).Value = ""   # CUL!L109 (cleared)
$ws.Cells.Item(109, 13).Value = 15000   # CUL!M109
$ws.Cells.Item(109, 14).Value = -17080   # CUL!N109
$ws.Cells.Item(110, 8).Value = 10499.5   # CUL!H110
$ws.Cells.Item(110, 10).Value = 13000   # CUL!J110
$ws.Cells.Item(110, 12).Value = 39000   # CUL!L110
$ws.Cells.Item(110, 14).Value = -47180   # CUL!N110
$ws.Cells.Item(115, 8).Value = 2850   # CUL!H115
$ws.Cells.Item(115, 9).Value = 2250   # CUL!I115
$ws.Cells.Item(115, 10).Value = 3450   # CUL!J115
$ws.Cells.Item(115, 11).Value = 6750   # CUL!K115
$ws.Cells.Item(115, 12).Value = 10350   # CUL!L115
$ws.Cells.Item(115, 13).Value = -5575   # CUL!M115
$ws.Cells.Item(115, 14).Value = -12700   # CUL!N115
$ws.Cells.Item(117, 8).Value = 1250   # CUL!H117
$ws.Cells.Item(117, 9).Value = 1300   # CUL!I117
$ws.Cells.Item(117, 10).Value = 1216.6666   # CUL!J117
$ws.Cells.Item(117, 11).Value = 3900   # CUL!K117
$ws.Cells.Item(117, 12).Value = 3649.9998   # CUL!L117
$ws.Cells.Item(117, 13).Value = -458   # CUL!M117
$ws.Cells.Item(117, 14).Value = -10533.9998   # CUL!N117
$ws.Cells.Item(118, 8).Value = 41669140   # CUL!H118
$ws.Cells.Item(118, 9).Value = 166667170   # CUL!I118
$ws.Cells.Item(118, 10).Value = 3130   # CUL!J118
$ws.Cells.Item(118, 11).Value = 500001510   # CUL!K118
$ws.Cells.Item(118, 12).Value = 9390   # CUL!L118
$ws.Cells.Item(118, 13).Value = -500000267   # CUL!M118
$ws.Cells.Item(118, 14).Value = -11876   # CUL!N118
$ws.Cells.Item(120, 8).Value = 576.6667   # CUL!H120
$ws.Cells.Item(120, 9).Value = 576.6667   # CUL!I120
$ws.Cells.Item(120, 10).Value = 0   # CUL!J120
$ws.Cells.Item(120, 11).Value = 1730.0001   # CUL!K120
$ws.Cells.Item(120, 12).Value = ""   # CUL!L120 (cleared)
$ws.Cells.Item(120, 13).Value = 3107.9999   # CUL!M120
$ws.Cells.Item(120, 14).Value = 0   # CUL!N120
$ws.Cells.Item(131, 8).Value = 602.23254   # CUL!H131
$ws.Cells.Item(131, 9).Value = 468.32257   # CUL!I131
$ws.Cells.Item(131, 10).Value = 948.1667   # CUL!J131
$ws.Cells.Item(131, 11).Value = 1404.96771   # CUL!K131
$ws.Cells.Item(131, 12).Value = 2844.5001   # CUL!L131
$ws.Cells.Item(131, 13).Value = 3635.03229   # CUL!M131
$ws.Cells.Item(131, 14).Value = -12924.5001   # CUL!N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 40429   # GSM!H138
$ws.Cells.Item(138, 10).Value = 40429   # GSM!J138
$ws.Cells.Item(138, 12).Value = 40429   # GSM!L138
$ws.Cells.Item(138, 14).Value = -50709   # GSM!N138

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1601.6086   # LTW!H61
$ws.Cells.Item(61, 9).Value = 1368.5714   # LTW!I61
$ws.Cells.Item(61, 11).Value = 1368.5714   # LTW!K61
$ws.Cells.Item(61, 13).Value = -1166.5714   # LTW!M61
$ws.Cells.Item(68, 8).Value = 2045   # LTW!H68
$ws.Cells.Item(68, 9).Value = 1950.7693   # LTW!I68
$ws.Cells.Item(68, 10).Value = 2220   # LTW!J68
$ws.Cells.Item(68, 11).Value = 1950.7693   # LTW!K68
$ws.Cells.Item(68, 12).Value = 2220   # LTW!L68
$ws.Cells.Item(68, 13).Value = -1201.7693   # LTW!M68
$ws.Cells.Item(68, 14).Value = -3718   # LTW!N68
$ws.Cells.Item(71, 8).Value = 2045   # LTW!H71
$ws.Cells.Item(71, 9).Value = 1950.7693   # LTW!I71
$ws.Cells.Item(71, 10).Value = 2220   # LTW!J71
$ws.Cells.Item(71, 11).Value = 9753.8465   # LTW!K71
$ws.Cells.Item(71, 12).Value = 11100   # LTW!L71
$ws.Cells.Item(71, 13).Value = -6009.8465   # LTW!M71
$ws.Cells.Item(71, 14).Value = -18588   # LTW!N71
$ws.Cells.Item(82, 8).Value = 824.619   # LTW!H82
$ws.Cells.Item(82, 9).Value = 554.5714   # LTW!I82
$ws.Cells.Item(82, 10).Value = 959.6429000000001   # LTW!J82
$ws.Cells.Item(82, 11).Value = 554.5714   # LTW!K82
$ws.Cells.Item(82, 12).Value = 959.6429000000001   # LTW!L82
$ws.Cells.Item(82, 13).Value = -193.5714   # LTW!M82
$ws.Cells.Item(82, 14).Value = -1681.6429   # LTW!N82
$ws.Cells.Item(85, 8).Value = 824.619   # LTW!H85
$ws.Cells.Item(85, 9).Value = 554.5714   # LTW!I85
$ws.Cells.Item(85, 10).Value = 959.6429000000001   # LTW!J85
$ws.Cells.Item(85, 11).Value = 554.5714   # LTW!K85
$ws.Cells.Item(85, 12).Value = 959.6429000000001   # LTW!L85
$ws.Cells.Item(85, 13).Value = 693.4286   # LTW!M85
$ws.Cells.Item(85, 14).Value = -3455.6429   # LTW!N85
$ws.Cells.Item(100, 8).Value = 1381.6471   # LTW!H100
$ws.Cells.Item(100, 9).Value = 1048   # LTW!I100
$ws.Cells.Item(100, 10).Value = 1858.2858   # LTW!J100
$ws.Cells.Item(100, 11).Value = 1048   # LTW!K100
$ws.Cells.Item(100, 12).Value = 1858.2858   # LTW!L100
$ws.Cells.Item(100, 13).Value = -507   # LTW!M100
$ws.Cells.Item(100, 14).Value = -2940.2858   # LTW!N100
$ws.Cells.Item(113, 8).Value = 1601.6086   # LTW!H113
$ws.Cells.Item(113, 9).Value = 1368.5714   # LTW!I113
$ws.Cells.Item(113, 11).Value = 1368.5714   # LTW!K113
$ws.Cells.Item(113, 13).Value = 801.4286   # LTW!M113

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 2648.5   # WVR!H62
$ws.Cells.Item(62, 9).Value = 2783.1428   # WVR!I62
$ws.Cells.Item(62, 10).Value = 2334.3333   # WVR!J62
$ws.Cells.Item(62, 11).Value = 2783.1428   # WVR!K62
$ws.Cells.Item(62, 12).Value = 2334.3333   # WVR!L62
$ws.Cells.Item(62, 13).Value = -2159.1428   # WVR!M62
$ws.Cells.Item(62, 14).Value = -3582.3333   # WVR!N62
$ws.Cells.Item(65, 8).Value = 2648.5   # WVR!H65
$ws.Cells.Item(65, 9).Value = 2783.1428   # WVR!I65
$ws.Cells.Item(65, 10).Value = 2334.3333   # WVR!J65
$ws.Cells.Item(65, 11).Value = 13915.714   # WVR!K65
$ws.Cells.Item(65, 12).Value = 11671.6665   # WVR!L65
$ws.Cells.Item(65, 13).Value = -10795.714   # WVR!M65
$ws.Cells.Item(65, 14).Value = -17911.6665   # WVR!N65
$ws.Cells.Item(140, 8).Value = 47485.8   # WVR!H140
$ws.Cells.Item(140, 10).Value = 47485.8   # WVR!J140
$ws.Cells.Item(140, 12).Value = 47485.8   # WVR!L140
$ws.Cells.Item(140, 14).Value = -57845.8   # WVR!N140
